$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

# Update the confidential disclosure date text in A10
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-28 for illustrative purposes only and are subject to change."

# Update holdings weight/percent-change values
$ws.Range("D2").Value = 0.4868203176890569
$ws.Range("E2").Value = -0.000386548125241748

$ws.Range("D3").Value = 0.3319931735417173
$ws.Range("E3").Value = 0.0007756447547022649

$ws.Range("D4").Value = 0.09716498945513993
$ws.Range("E4").Value = -0.003389225829468545

$ws.Range("D5").Value = 0.05398398864980603
$ws.Range("E5").Value = 0.001949317738791478

$ws.Range("D6").Value = 0.03003753066427976
$ws.Range("E6").Value = 0.02751159196290587

$ws.Range("E7").Value = 0.0006716274244356057

# Restore worksheet protection (cannot reproduce the original legacy
# password hash, but re-apply protection so the sheet isn't left open)
$ws.Protect()
